# Updated cryptos list on Tue May  7 08:41:56 UTC 2024 with GitHub Actions
#
# Refreshes the Coin / Link / Price / Volume(1h) columns (B:E) on rows
# 2-51 of the crypto-ranking sheet with the latest scraped snapshot.
# Two pairs of coins (rows 15/16 and 45/46) swapped rank order between
# scrapes, so those rows are rewritten in full (Coin, Link, Price,
# Volume); the rest only get Price/Volume refreshed.
#
# Price strings that happen to look like a plain number (e.g. "0.120",
# "1.00", "0.0000240") are written with a temporary Text number format
# so Excel's automatic type detection doesn't normalize away the
# significant trailing/leading zeros; the format is cleared again right
# after the write so the cell is left with no explicit style, same as
# before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Range, $Value) {
    if ($Value -match '^[+-]?\d+(\.\d+)?$') {
        $Range.NumberFormat = "@"
        $Range.Value = $Value
        $Range.ClearFormats()
    } else {
        $Range.Value = $Value
    }
}


# Row 2
Set-TextValue $ws.Range("D2") '64.251.57'
Set-TextValue $ws.Range("E2") '  -2.03%  '

# Row 3
Set-TextValue $ws.Range("D3") '3.114.26'
Set-TextValue $ws.Range("E3") '  -3.08%  '

# Row 4
Set-TextValue $ws.Range("E4") '  -0.10%  '

# Row 5
Set-TextValue $ws.Range("D5") '592.24'
Set-TextValue $ws.Range("E5") '  -1.46%  '

# Row 6
Set-TextValue $ws.Range("D6") '157.05'
Set-TextValue $ws.Range("E6") '  +2.83%  '

# Row 7
Set-TextValue $ws.Range("E7") '  -0.07%  '

# Row 8
Set-TextValue $ws.Range("E8") '  -1.01%  '

# Row 9
Set-TextValue $ws.Range("D9") '3.113.85'
Set-TextValue $ws.Range("E9") '  -2.93%  '

# Row 10
Set-TextValue $ws.Range("E10") '  -5.12%  '

# Row 11
Set-TextValue $ws.Range("D11") '5.93'
Set-TextValue $ws.Range("E11") '  -4.28%  '

# Row 12
Set-TextValue $ws.Range("E12") '  -4.25%  '

# Row 13
Set-TextValue $ws.Range("D13") '37.24'
Set-TextValue $ws.Range("E13") '  -5.76%  '

# Row 14
Set-TextValue $ws.Range("D14") '0.0000240'
Set-TextValue $ws.Range("E14") '  -5.78%  '

# Row 15
Set-TextValue $ws.Range("B15") 'WrappedliquidstakedEther2.0'
Set-TextValue $ws.Range("C15") 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws.Range("D15") '3.621.66'
Set-TextValue $ws.Range("E15") '  -3.36%  '

# Row 16
Set-TextValue $ws.Range("B16") 'TRON'
Set-TextValue $ws.Range("C16") 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Range("D16") '0.120'
Set-TextValue $ws.Range("E16") '  -1.67%  '

# Row 17
Set-TextValue $ws.Range("D17") '7.24'
Set-TextValue $ws.Range("E17") '  -2.50%  '

# Row 18
Set-TextValue $ws.Range("D18") '64.076.42'
Set-TextValue $ws.Range("E18") '  -1.72%  '

# Row 19
Set-TextValue $ws.Range("D19") '3.110.95'
Set-TextValue $ws.Range("E19") '  -6.62%  '

# Row 20
Set-TextValue $ws.Range("D20") '480.96'
Set-TextValue $ws.Range("E20") '  -0.58%  '

# Row 21
Set-TextValue $ws.Range("D21") '14.52'
Set-TextValue $ws.Range("E21") '  -3.37%  '

# Row 22
Set-TextValue $ws.Range("E22") '  -7.73%  '

# Row 23
Set-TextValue $ws.Range("E23") '  -3.94%  '

# Row 24
Set-TextValue $ws.Range("D24") '2.45'
Set-TextValue $ws.Range("E24") '  -2.37%  '

# Row 25
Set-TextValue $ws.Range("D25") '12.96'
Set-TextValue $ws.Range("E25") '  -5.99%  '

# Row 26
Set-TextValue $ws.Range("D26") '81.37'
Set-TextValue $ws.Range("E26") '  -2.78%  '

# Row 27
Set-TextValue $ws.Range("D27") '10.44'
Set-TextValue $ws.Range("E27") '  +3.78%  '

# Row 28
Set-TextValue $ws.Range("E28") '  -0.32%  '

# Row 29
Set-TextValue $ws.Range("D29") '7.48'
Set-TextValue $ws.Range("E29") '  -0.86%  '

# Row 30
Set-TextValue $ws.Range("E30") '  -4.13%  '

# Row 31
Set-TextValue $ws.Range("D31") '1.00'
Set-TextValue $ws.Range("E31") '  -0.03%  '

# Row 32
Set-TextValue $ws.Range("E32") '  -4.12%  '

# Row 33
Set-TextValue $ws.Range("D33") '0.114'
Set-TextValue $ws.Range("E33") '  -6.26%  '

# Row 34
Set-TextValue $ws.Range("D34") '27.45'
Set-TextValue $ws.Range("E34") '  -4.91%  '

# Row 35
Set-TextValue $ws.Range("D35") '0.0₃0842'
Set-TextValue $ws.Range("E35") '  -5.13%  '

# Row 36
Set-TextValue $ws.Range("E36") '  -2.63%  '

# Row 37
Set-TextValue $ws.Range("E37") '  -5.20%  '

# Row 38
Set-TextValue $ws.Range("D38") '3.29'
Set-TextValue $ws.Range("E38") '  -6.64%  '

# Row 39
Set-TextValue $ws.Range("D39") '2.25'
Set-TextValue $ws.Range("E39") '  -6.61%  '

# Row 40
Set-TextValue $ws.Range("D40") '51.11'
Set-TextValue $ws.Range("E40") '  -2.35%  '

# Row 41
Set-TextValue $ws.Range("D41") '9.21'
Set-TextValue $ws.Range("E41") '  -2.86%  '

# Row 42
Set-TextValue $ws.Range("D42") '439.83'
Set-TextValue $ws.Range("E42") '  -9.40%  '

# Row 43
Set-TextValue $ws.Range("D43") '0.291'
Set-TextValue $ws.Range("E43") '  -4.48%  '

# Row 44
Set-TextValue $ws.Range("E44") '  -5.53%  '

# Row 45
Set-TextValue $ws.Range("B45") 'Kaspa'
Set-TextValue $ws.Range("C45") 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range("D45") '0.112'
Set-TextValue $ws.Range("E45") '  -0.59%  '

# Row 46
Set-TextValue $ws.Range("B46") 'Arweave'
Set-TextValue $ws.Range("C46") 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-TextValue $ws.Range("D46") '40.16'
Set-TextValue $ws.Range("E46") '  +1.78%  '

# Row 47
Set-TextValue $ws.Range("D47") '2.833.94'
Set-TextValue $ws.Range("E47") '  -4.00%  '

# Row 48
Set-TextValue $ws.Range("D48") '130.16'
Set-TextValue $ws.Range("E48") '  -1.51%  '

# Row 49
Set-TextValue $ws.Range("D49") '25.39'
Set-TextValue $ws.Range("E49") '  -0.35%  '

# Row 50
Set-TextValue $ws.Range("E50") '  +0.04%  '

# Row 51
Set-TextValue $ws.Range("E51") '  -3.77%  '
